$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look numeric,
# so Excel stores them as text (matching the source data) instead of
# auto-converting to a number and losing formatting (e.g. trailing zeros).
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '43.676.17'
$ws.Range('E2').Value = '  +2.99%  '
$ws.Range('D3').Value = '2.201.25'
$ws.Range('E3').Value = '  +0.68%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '258.86'
$ws.Range('E5').Value = '  +2.45%  '
$ws.Range('D6').Value = '83.10'
$ws.Range('E6').Value = '  +10.60%  '
$ws.Range('D7').Value = '0.620'
$ws.Range('E7').Value = '  +1.05%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +2.41%  '
$ws.Range('D10').Value = '44.34'
$ws.Range('E10').Value = '  +9.68%  '
$ws.Range('D11').Value = '0.0914'
$ws.Range('E11').Value = '  +0.33%  '
$ws.Range('D12').Value = '7.17'
$ws.Range('E12').Value = '  +5.66%  '
$ws.Range('E13').Value = '  +2.45%  '
$ws.Range('D14').Value = '2.527.75'
$ws.Range('E14').Value = '  +0.60%  '
$ws.Range('D15').Value = '14.36'
$ws.Range('E15').Value = '  +1.68%  '
$ws.Range('D16').Value = '2.194.25'
$ws.Range('E16').Value = '  +0.66%  '
$ws.Range('D17').Value = '0.780'
$ws.Range('E17').Value = '  +1.18%  '
$ws.Range('D18').Value = '43.595.86'
$ws.Range('E18').Value = '  +2.96%  '
$ws.Range('E19').Value = '  +1.80%  '
$ws.Range('D20').Value = '69.69'
$ws.Range('E21').Value = '  +0.72%  '
$ws.Range('E22').Value = '  +11.40%  '
$ws.Range('D23').Value = '231.49'
$ws.Range('E23').Value = '  +2.53%  '
$ws.Range('D24').Value = '8.96'
$ws.Range('E24').Value = '  -4.72%  '
$ws.Range('D26').Value = '10.60'
$ws.Range('E26').Value = '  +1.40%  '
$ws.Range('D27').Value = '3.46'
$ws.Range('E27').Value = '  +2.64%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '2.27'
$ws.Range('E28').Value = '  +3.16%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').Value = '2.24'
$ws.Range('E29').Value = '  +3.52%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').Value = '38.98'
$ws.Range('E30').Value = '  +1.36%  '
$ws.Range('D31').Value = '174.16'
$ws.Range('E31').Value = '  +0.88%  '
$ws.Range('D32').Value = '20.39'
$ws.Range('E32').Value = '  +1.73%  '
$ws.Range('E33').Value = '  +4.66%  '
$ws.Range('D34').Value = '5.31'
$ws.Range('E34').Value = '  +3.63%  '
$ws.Range('E35').Value = '  +2.09%  '
$ws.Range('E36').Value = '  +3.90%  '
$ws.Range('D37').Value = '4.51'
$ws.Range('E37').Value = '  +6.62%  '
$ws.Range('D38').Value = '0.0358'
$ws.Range('E38').Value = '  +5.87%  '
$ws.Range('D39').Value = '12.49'
$ws.Range('E39').Value = '  +3.98%  '
$ws.Range('E40').Value = '  +8.00%  '
$ws.Range('E41').Value = '  +1.64%  '
$ws.Range('D42').Value = '62.85'
$ws.Range('E42').Value = '  +6.80%  '
$ws.Range('D43').Value = '5.48'
$ws.Range('E43').Value = '  +6.17%  '
$ws.Range('E44').Value = '  +3.22%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = '8.35'
$ws.Range('E45').Value = '  +2.32%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').Value = '0.0977'
$ws.Range('E46').Value = '  +0.75%  '
$ws.Range('D47').Value = '99.66'
$ws.Range('E47').Value = '  -1.67%  '
$ws.Range('E48').Value = '  +5.80%  '
$ws.Range('E49').Value = '  +1.86%  '
$ws.Range('D50').Value = '0.438'
$ws.Range('E50').Value = '  -4.67%  '
$ws.Range('D51').Value = '1.48'
$ws.Range('E51').Value = '  +7.05%  '
